# Generate Report for Handoff
# Regenerates the localization-status report: the existing "Ready for
# handoff" row is replaced with a fresh CI run (new source UUIDs / new
# handoff timestamps) and two more rows are appended for the extra files
# that are now part of the handoff set.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item(1)   # "Overview"
$zhcn     = $wb.Worksheets.Item(2)   # "zh-cn"
$dede     = $wb.Worksheets.Item(3)   # "de-de"

# The engine's Range.Hyperlinks.Delete() removes every hyperlink on the
# sheet (not just the range), so clear each sheet's hyperlinks exactly once
# up front and rebuild them all from scratch below.
$overview.Hyperlinks.Delete()
$zhcn.Hyperlinks.Delete()
$dede.Hyperlinks.Delete()

$baseRepo   = "https://github.com/OpenLocalizationTest/oltest/blob/0c249cacbe3fb58c2e0e896f611ba180cd851d0d/e2e/"
$zhHandback = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0e165907aa40f4da7e6b4dec5c8c4c38fe9e94e6/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/"
$deHandback = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/b804ff3e44c24c0e08f6a540be97a391e6c4433f/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/"

$file1 = "db9e7ad0-c1df-4b46-9f29-116537b786b1.png"
$file2 = "e18de121-bf5a-4a99-a111-41dc10cbcddb.md"
$file3 = "facb285f-69a3-4797-97ae-f4c988294210.png"

$xlf1zh = "f3596fd849ef0addca00f88e634a43d82d7e43bb.png"
$xlf2zh = "e18de121-bf5a-4a99-a111-41dc10cbcddb.27f715458cab92f044f809e282f757658ab27e82.zh-cn.xlf"
$xlf3zh = "01d4a395a978a59c806237006d4269deca95afa1.png"

$xlf1de = "f3596fd849ef0addca00f88e634a43d82d7e43bb.png"
$xlf2de = "e18de121-bf5a-4a99-a111-41dc10cbcddb.27f715458cab92f044f809e282f757658ab27e82.de-de.xlf"
$xlf3de = "01d4a395a978a59c806237006d4269deca95afa1.png"

$overviewDate = "2016-25-17 10:25:43"
$handoffDate  = "2016-03-17 10:25:40"
$handoffDateZh = "2016-03-17 10:25:40"
$handoffDateDe = "2016-03-17 10:25:43"
$zeroDate     = "0001-01-01 00:00:00"
$status       = "Ready for handoff"

function Set-Overview-Row($ws, [int]$row, [string]$fileName, [string]$linkUrl) {
    $ws.Range("A$row").Value = $fileName
    $ws.Hyperlinks.Add($ws.Range("A$row"), $linkUrl, [Type]::Missing, [Type]::Missing, $fileName) | Out-Null
    $ws.Range("A$row").Style = "HyperLink"

    $ws.Range("B$row").Value = $status
    $ws.Range("C$row").Value = $status
    $ws.Range("D$row").Value = $overviewDate
}

$file1Url = $baseRepo + $file1
$file2Url = $baseRepo + $file2
$file3Url = $baseRepo + $file3

Set-Overview-Row $overview 2 $file1 $file1Url
Set-Overview-Row $overview 3 $file2 $file2Url
Set-Overview-Row $overview 4 $file3 $file3Url

function Set-Lang-Row($ws, [int]$row, [string]$sourceFile, [string]$ext, [string]$targetFile, [string]$handoffDisplayDate, [string]$handbackBase, [string]$reason, [string]$dependencyFrom) {
    $sourceUrl = $baseRepo + $sourceFile
    $targetUrl = $handbackBase + $sourceFile + "." + $targetFile

    $ws.Range("A$row").Value = $sourceFile
    $ws.Hyperlinks.Add($ws.Range("A$row"), $sourceUrl, [Type]::Missing, [Type]::Missing, $sourceFile) | Out-Null
    $ws.Range("A$row").Style = "HyperLink"

    $ws.Range("B$row").Value = $ext
    $ws.Hyperlinks.Add($ws.Range("B$row"), $sourceUrl, [Type]::Missing, [Type]::Missing, $ext) | Out-Null
    $ws.Range("B$row").Style = "HyperLink"

    $ws.Range("C$row").Value = $status

    $ws.Range("D$row").Value = $targetFile
    $ws.Hyperlinks.Add($ws.Range("D$row"), $targetUrl, [Type]::Missing, [Type]::Missing, $targetFile) | Out-Null
    $ws.Range("D$row").Style = "HyperLink"

    $ws.Range("E$row").Value = $handoffDisplayDate
    $ws.Range("E$row").NumberFormat = "yyyy-mm-dd HH:mm:ss"

    $ws.Range("H$row").Value = $zeroDate

    $ws.Range("I$row").Value = $reason

    if ($dependencyFrom) {
        $ws.Range("J$row").Value = $dependencyFrom
    }
}

$dependencyFromFile2 = "e2e\" + $file2

# zh-cn sheet
Set-Lang-Row $zhcn 2 $file1 ".png" $xlf1zh $handoffDateZh $zhHandback "IsDependency" $dependencyFromFile2
Set-Lang-Row $zhcn 3 $file2 ".md"  $xlf2zh $handoffDateZh $zhHandback "Include" $null
Set-Lang-Row $zhcn 4 $file3 ".png" $xlf3zh $handoffDateZh $zhHandback "IsDependency" $dependencyFromFile2

# de-de sheet
Set-Lang-Row $dede 2 $file1 ".png" $xlf1de $handoffDateDe $deHandback "IsDependency" $dependencyFromFile2
Set-Lang-Row $dede 3 $file2 ".md"  $xlf2de $handoffDateDe $deHandback "Include" $null
Set-Lang-Row $dede 4 $file3 ".png" $xlf3de $handoffDateDe $deHandback "IsDependency" $dependencyFromFile2
